$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.048.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.421.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.49%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.59%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9975"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "276.51"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3699"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.93%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3137"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.66"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.057"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06537"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9988"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.513"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.82"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.203"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.422.25"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001022"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05695"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9981"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.66"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.616"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.87"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.10"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.231"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.099.71"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.288"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.04"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.35"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.582.33"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.96"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.945"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +9.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.271"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8255"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -9.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07797"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.471"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.934"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05849"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.039"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9973"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.60"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.68%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.32%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1878"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5335"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.41"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.546"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "117.72"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5226"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.782"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9958"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.96%  "

# Reset number format/style back to default (no explicit style) for touched cells
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Style = "Normal"
